$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the direct formatting of an existing, "typical" data row (row 65 uses
# the standard style set: A=2, B=7, C=11, D=2, E=7) down onto the new row 70
# so the new row picks up identical borders/fonts/alignment without touching
# styles.xml.
$ws.Range("A65:E65").Copy()
$ws.Range("A70:E70").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new test case data.
$ws.Range("A70").Value = "Profile69"
$ws.Range("B70").Value = "OPQA-500"
$ws.Range("C70").Value = "Verify that user is able to update his own profile picture"
$ws.Range("D70").Value = "Y"
$ws.Range("E70").Value = ""

# Match the updated view state captured in the saved workbook (best effort —
# the headless view model only persists the selection/active cell).
$excel.ActiveWindow.ScrollRow = 56
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("F67").Select()
